$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at C, shifting compliance/ID/mod to D/E/F
$ws.Columns("C").Insert()

# Header
$ws.Range("C1").Value = "unique_roles"

# Populate unique_roles values for rows 2-28
$ws.Range("C2").Value = "Dairy Cow,Milk Mage"
$ws.Range("C3").Value = "Ranger"
$ws.Range("C4").Value = "Mule,Siege,Wool Mage"
$ws.Range("C5").Value = "Ice Mage,Lancer,Love Mage"
$ws.Range("C6").Value = "Vine Mage,Wood Mage"
$ws.Range("C7").Value = "Hypnomancer,Lancer"
$ws.Range("C8").Value = "Ice Mage"
$ws.Range("C9").Value = "Assassin,Waitress"
$ws.Range("C10").Value = "Lancer,Silk Mage"
$ws.Range("C11").Value = "Hunter,Lancer,Mule"
$ws.Range("C12").Value = "Ranger,Sharpshooter"
$ws.Range("C13").Value = "Feral,Pet,Rogue"
$ws.Range("C14").Value = "Songmage"
$ws.Range("C15").Value = "Barbarian"
$ws.Range("C16").Value = "Druid,Feral,Pet"
$ws.Range("C17").Value = "Hunter,Rogue"
$ws.Range("C18").Value = "Surgeon"
$ws.Range("C19").Value = "Vine Mage,Wood Mage"
$ws.Range("C20").Value = "Electro Mage"
$ws.Range("C21").Value = "Arousal Mage,Latex Mage,Metal Mage,Milk Mage"
$ws.Range("C22").Value = "Sculptor"
$ws.Range("C23").Value = "Surgeon"
$ws.Range("C24").Value = "Drain Mage,Love Mage"
$ws.Range("C25").Value = "Heavy Warrior,Shaved,Smith"
$ws.Range("C26").Value = "Leaf Dancer,Pocket Healer,Prankster"
$ws.Range("C27").Value = "Breeder,Improviser,Scoundrel"
$ws.Range("C28").Value = "Ambusher,Opportunist,Underling"
